$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Add()
$ws.Name = "Scatter"

$aVals = @(0,1,2,3,4,5,6,7,8,9)
$bVals = @(0,1,4,9,16,25,36,49,64,81)
for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item($i+1, 1).Value = $aVals[$i]
    $ws.Cells.Item($i+1, 2).Value = $bVals[$i]
}

$chartObj = $ws.Shapes.AddChart2(-1, -4169)
$chart = $chartObj.Chart
$chart.SetSourceData($ws.Range("A1:B10"))

$ws.Move($null, $wb.Worksheets.Item("Line"))
